# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list ... with GitHub Actions" — only column D (Price) and
# column E (Volume(1h)) text for rows 2-51 change; everything else (coin name,
# link, row styling) is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.915.52'
$ws.Range("E2").Value = '  -0.12%  '

$ws.Range("D3").Value = '1.635.93'
$ws.Range("E3").Value = '  -0.33%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.22'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5069'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.00%  '

$ws.Range("E7").Value = '  -0.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2579'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.64%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06365'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.12%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.62'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.68%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07755'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.24%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.264'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.55%  '

$ws.Range("D13").Value = '1.629.23'
$ws.Range("E13").Value = '  -1.11%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5528'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.42%  '

$ws.Range("D15").Value = '0.0₅7732'
$ws.Range("E15").Value = '  -1.27%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.15'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.21%  '

$ws.Range("D17").Value = '25.921.46'
$ws.Range("E17").Value = '  -0.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.002'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.443'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '194.59'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.47%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.904'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.062'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.21%  '

$ws.Range("E23").Value = '  -0.19%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.907'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.44%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '142.85'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1241'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +6.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.837'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.66%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.58'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.82%  '

$ws.Range("E29").Value = '  +0.60%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.04869'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.253'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.31%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.194'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.30%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.547'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.370'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.36%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9062'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.30%  '

$ws.Range("E36").Value = '  -0.78%  '

$ws.Range("E37").Value = '  +1.21%  '

$ws.Range("D38").Value = '1.122.53'
$ws.Range("E38").Value = '  -0.92%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01559'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.30%  '

$ws.Range("E40").Value = '  -0.26%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.579'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.36%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8056'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.25%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '97.69'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.06%  '

$ws.Range("D44").Value = '0.0₈120'
$ws.Range("E44").Value = '  -6.04%  '

$ws.Range("D45").Value = '1.773.13'
$ws.Range("E45").Value = '  -0.25%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4458'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.75%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.97'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.35%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9968'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.67%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05153'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.59%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.561'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.67%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.004'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.04%  '
